$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the numeric-looking columns (D, E, G) so values
# are stored as literal strings instead of being auto-converted to numbers,
# matching the original inlineStr cell type.
$textRangeDE = $ws.Range("D2:E51")
$textRangeG = $ws.Range("G2:G51")
$textRangeDE.NumberFormat = "@"
$textRangeG.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '310.30'
$ws.Range("E2").Value = '-2.61%'
$ws.Range("G2").Value = '15'
# Row 3
$ws.Range("D3").Value = '54.26'
$ws.Range("E3").Value = '12.58%'
$ws.Range("G3").Value = '15'
# Row 4
$ws.Range("D4").Value = '5.129'
$ws.Range("E4").Value = '-2.65%'
$ws.Range("G4").Value = '15'
# Row 5
$ws.Range("D5").Value = '0.07801'
$ws.Range("E5").Value = '-1.57%'
$ws.Range("G5").Value = '15'
# Row 6
$ws.Range("D6").Value = '4.518'
$ws.Range("E6").Value = '-1.41%'
$ws.Range("G6").Value = '15'
# Row 7
$ws.Range("D7").Value = '1.359'
$ws.Range("E7").Value = '0.24%'
$ws.Range("G7").Value = '15'
# Row 8
$ws.Range("D8").Value = '1.577'
$ws.Range("E8").Value = '-3.33%'
$ws.Range("G8").Value = '15'
# Row 9
$ws.Range("D9").Value = '0.1214'
$ws.Range("E9").Value = '-5.93%'
$ws.Range("G9").Value = '15'
# Row 10
$ws.Range("E10").Value = '2.08%'
$ws.Range("G10").Value = '15'
# Row 11
$ws.Range("D11").Value = '0.04730'
$ws.Range("E11").Value = '2.53%'
$ws.Range("G11").Value = '15'
# Row 12
$ws.Range("D12").Value = '0.09433'
$ws.Range("E12").Value = '0.72%'
$ws.Range("G12").Value = '15'
# Row 13
$ws.Range("D13").Value = '0.1043'
$ws.Range("E13").Value = '-0.50%'
$ws.Range("G13").Value = '15'
# Row 14
$ws.Range("D14").Value = '0.001262'
$ws.Range("E14").Value = '-4.72%'
$ws.Range("G14").Value = '15'
# Row 15
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").Value = '0.04177'
$ws.Range("E15").Value = '0.33%'
$ws.Range("G15").Value = '15'
# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005803'
$ws.Range("E16").Value = '-2.27%'
$ws.Range("G16").Value = '15'
# Row 17
$ws.Range("B17").Value = 'UpBots'
$ws.Range("C17").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D17").Value = '0.007487'
$ws.Range("E17").Value = '2,013.54%'
$ws.Range("G17").Value = '15'
# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '3.342'
$ws.Range("E18").Value = '0.10%'
$ws.Range("G18").Value = '15'
# Row 19
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '2.414'
$ws.Range("E19").Value = '-0.54%'
$ws.Range("G19").Value = '15'
# Row 20
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '0.3431'
$ws.Range("E20").Value = '-0.79%'
$ws.Range("G20").Value = '15'
# Row 21
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").Value = '8.041'
$ws.Range("E21").Value = '-1.29%'
$ws.Range("G21").Value = '15'
# Row 22
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = '0.1363'
$ws.Range("E22").Value = '-1.37%'
$ws.Range("G22").Value = '15'
# Row 23
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '0.3093'
$ws.Range("E23").Value = '-0.13%'
$ws.Range("G23").Value = '15'
# Row 24
$ws.Range("D24").Value = '0.001259'
$ws.Range("E24").Value = '-4.56%'
$ws.Range("G24").Value = '15'
# Row 25
$ws.Range("D25").Value = '0.003912'
$ws.Range("E25").Value = '-8.08%'
$ws.Range("G25").Value = '15'
# Row 26
$ws.Range("D26").Value = '0.0001350'
$ws.Range("E26").Value = '-0.08%'
$ws.Range("G26").Value = '15'
# Row 27
$ws.Range("G27").Value = '15'
# Row 28
$ws.Range("G28").Value = '15'
# Row 29
$ws.Range("G29").Value = '15'
# Row 30
$ws.Range("G30").Value = '15'
# Row 31
$ws.Range("G31").Value = '15'
# Row 32
$ws.Range("G32").Value = '15'
# Row 33
$ws.Range("G33").Value = '15'
# Row 34
$ws.Range("G34").Value = '15'
# Row 35
$ws.Range("G35").Value = '15'
# Row 36
$ws.Range("G36").Value = '15'
# Row 37
$ws.Range("G37").Value = '15'
# Row 38
$ws.Range("D38").Value = '0.02589'
$ws.Range("E38").Value = '-3.01%'
$ws.Range("G38").Value = '15'
# Row 39
$ws.Range("D39").Value = '0.05914'
$ws.Range("E39").Value = '3.61%'
$ws.Range("G39").Value = '15'
# Row 40
$ws.Range("D40").Value = '0.01053'
$ws.Range("E40").Value = '-2.12%'
$ws.Range("G40").Value = '15'
# Row 41
$ws.Range("D41").Value = '0.007892'
$ws.Range("E41").Value = '-1.59%'
$ws.Range("G41").Value = '15'
# Row 42
$ws.Range("D42").Value = '0.1417'
$ws.Range("E42").Value = '-1.21%'
$ws.Range("G42").Value = '15'
# Row 43
$ws.Range("D43").Value = '0.008226'
$ws.Range("E43").Value = '10.51%'
$ws.Range("G43").Value = '15'
# Row 44
$ws.Range("D44").Value = '0.008484'
$ws.Range("E44").Value = '0.19%'
$ws.Range("G44").Value = '15'
# Row 45
$ws.Range("D45").Value = '0.3111'
$ws.Range("E45").Value = '-2.10%'
$ws.Range("G45").Value = '15'
# Row 46
$ws.Range("D46").Value = '0.00007280'
$ws.Range("E46").Value = '9.93%'
$ws.Range("G46").Value = '15'
# Row 47
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").Value = '-0.10%'
$ws.Range("G47").Value = '15'
# Row 48
$ws.Range("D48").Value = '0.05663'
$ws.Range("E48").Value = '3.04%'
$ws.Range("G48").Value = '15'
# Row 49
$ws.Range("D49").Value = '0.002619'
$ws.Range("E49").Value = '-34.58%'
$ws.Range("G49").Value = '15'
# Row 50
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").Value = '-0.10%'
$ws.Range("G50").Value = '15'
# Row 51
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").Value = '-0.10%'
$ws.Range("G51").Value = '15'

# Restore default style (removes the Text number-format override we applied
# above) while keeping the values stored as text.
$textRangeDE.Style = "Normal"
$textRangeG.Style = "Normal"
